{"js": "// Replace each \"AxB=C\" answer text in the table with the new value.\n// Each old value is unique in the document, so a plain exact-match\n// search + full-match replace is reliable and keeps run formatting\n// (rFonts/sz) intact because insertText(\"Replace\") only swaps the text\n// inside the matched range.\nconst replacements = [\n  [\"62\u00d758=3596\", \"55\u00d766=3630\"],\n  [\"22\u00d739=858\", \"95\u00d786=8170\"],\n  [\"79\u00d773=5767\", \"68\u00d793=6324\"],\n  [\"39\u00d724=936\", \"54\u00d751=2754\"],\n  [\"74\u00d734=2516\", \"30\u00d732=960\"],\n  [\"43\u00d721=903\", \"21\u00d749=1029\"],\n  [\"45\u00d716=720\", \"81\u00d766=5346\"],\n  [\"54\u00d771=3834\", \"72\u00d794=6768\"],\n  [\"40\u00d713=520\", \"27\u00d738=1026\"],\n  [\"35\u00d733=1155\", \"94\u00d755=5170\"],\n  [\"94\u00d753=4982\", \"82\u00d711=902\"],\n  [\"65\u00d725=1625\", \"31\u00d736=1116\"],\n  [\"50\u00d790=4500\", \"58\u00d774=4292\"],\n  [\"17\u00d765=1105\", \"56\u00d774=4144\"],\n  [\"31\u00d725=775\", \"41\u00d733=1353\"],\n  [\"18\u00d753=954\", \"29\u00d783=2407\"],\n  [\"59\u00d780=4720\", \"64\u00d794=6016\"],\n  [\"46\u00d732=1472\", \"65\u00d740=2600\"],\n  [\"21\u00d727=567\", \"44\u00d756=2464\"],\n  [\"71\u00d793=6603\", \"42\u00d768=2856\"],\n  [\"64\u00d772=4608\", \"66\u00d798=6468\"],\n  [\"44\u00d723=1012\", \"58\u00d732=1856\"],\n  [\"25\u00d799=2475\", \"77\u00d773=5621\"],\n  [\"81\u00d792=7452\", \"76\u00d793=7068\"],\n  [\"84\u00d723=1932\", \"69\u00d727=1863\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"AxB=C\" answer text in the table with the new value,\n# using Find/Replace on the whole document content range. Each old\n# value is unique in the document, so a MatchWholeWord-style exact\n# literal Find/Replace (MatchWildcards off) swaps only the text run\n# content and leaves run formatting (rFonts/sz) untouched.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll    = 2\n$wdFindContinue  = 1\n\n$replacements = @(\n  @{ Old = \"62\u00d758=3596\"; New = \"55\u00d766=3630\" },\n  @{ Old = \"22\u00d739=858\"; New = \"95\u00d786=8170\" },\n  @{ Old = \"79\u00d773=5767\"; New = \"68\u00d793=6324\" },\n  @{ Old = \"39\u00d724=936\"; New = \"54\u00d751=2754\" },\n  @{ Old = \"74\u00d734=2516\"; New = \"30\u00d732=960\" },\n  @{ Old = \"43\u00d721=903\"; New = \"21\u00d749=1029\" },\n  @{ Old = \"45\u00d716=720\"; New = \"81\u00d766=5346\" },\n  @{ Old = \"54\u00d771=3834\"; New = \"72\u00d794=6768\" },\n  @{ Old = \"40\u00d713=520\"; New = \"27\u00d738=1026\" },\n  @{ Old = \"35\u00d733=1155\"; New = \"94\u00d755=5170\" },\n  @{ Old = \"94\u00d753=4982\"; New = \"82\u00d711=902\" },\n  @{ Old = \"65\u00d725=1625\"; New = \"31\u00d736=1116\" },\n  @{ Old = \"50\u00d790=4500\"; New = \"58\u00d774=4292\" },\n  @{ Old = \"17\u00d765=1105\"; New = \"56\u00d774=4144\" },\n  @{ Old = \"31\u00d725=775\"; New = \"41\u00d733=1353\" },\n  @{ Old = \"18\u00d753=954\"; New = \"29\u00d783=2407\" },\n  @{ Old = \"59\u00d780=4720\"; New = \"64\u00d794=6016\" },\n  @{ Old = \"46\u00d732=1472\"; New = \"65\u00d740=2600\" },\n  @{ Old = \"21\u00d727=567\"; New = \"44\u00d756=2464\" },\n  @{ Old = \"71\u00d793=6603\"; New = \"42\u00d768=2856\" },\n  @{ Old = \"64\u00d772=4608\"; New = \"66\u00d798=6468\" },\n  @{ Old = \"44\u00d723=1012\"; New = \"58\u00d732=1856\" },\n  @{ Old = \"25\u00d799=2475\"; New = \"77\u00d773=5621\" },\n  @{ Old = \"81\u00d792=7452\"; New = \"76\u00d793=7068\" },\n  @{ Old = \"84\u00d723=1932\"; New = \"69\u00d727=1863\" }\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair.Old\n  $find.Replacement.Text = $pair.New\n  $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair.New, $wdReplaceAll) | Out-Null\n}\n"}
